$d = $word.ActiveDocument

$replacements = @(
    @("30÷7=", "21÷8="),
    @("50÷5=", "53÷6="),
    @("29÷8=", "81÷2="),
    @("95÷9=", "83÷4="),
    @("77÷5=", "14÷4="),
    @("99÷5=", "95÷3="),
    @("10÷8=", "85÷9="),
    @("10÷2=", "20÷9="),
    @("45÷5=", "12÷7="),
    @("17÷4=", "98÷8="),
    @("24÷7=", "31÷3="),
    @("24÷8=", "19÷7="),
    @("31÷5=", "63÷4="),
    @("24÷2=", "25÷2="),
    @("14÷5=", "91÷6="),
    @("70÷3=", "65÷8="),
    @("52÷9=", "15÷9="),
    @("76÷3=", "25÷2="),
    @("91÷8=", "13÷9="),
    @("83÷3=", "21÷3="),
    @("77÷3=", "36÷8="),
    @("64÷7=", "97÷2="),
    @("74÷6=", "88÷2="),
    @("57÷9=", "66÷8="),
    @("17÷5=", "12÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
